# Add columns I ("I0") and J ("IF") to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row (row 1): copy style from H1 (bold/border/centered header style)
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:J19
$data = @(
    @(4, 6),
    @(9, 9),
    @(5, 6),
    @(1, 4),
    @(1, 6),
    @(1, 7),
    @(1, 6),
    @(1, 4),
    @(1, 4),
    @(1, 2),
    @(1, 3),
    @(1, 3),
    @(1, 5),
    @(1, 5),
    @(1, 5),
    @(1, 4),
    @(4, 6),
    @(4, 5)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]   # column I
    $ws.Cells.Item($row, 10).Value = $pair[1]  # column J
    $row++
}
